$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Empresa")

# Rename existing company name from full name to short acronym "UMG"
$ws.Range("B2").Value = "UMG"

# Append a new row (row 8) with the new company data
$ws.Range("A8").Value = '"67c1adb6a7cf6af9db40f795"'
$ws.Range("B8").Value = "HyperX"
$ws.Range("C8").Value = "MEDIO"
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = "Joaquin Figueroa, Saul de Leon"
$ws.Range("F8").Value = "Distribuidora"
$ws.Range("G8").Value = $true
